$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (O) by copying formatting from column N (2020)
# and filling in the new values, row by row, matching the source diff.

$values = @{
    3  = 2021
    5  = 2148.2
    6  = 109.5
    7  = 210.1
    8  = 196
    9  = 209
    10 = 300.2
    11 = 302.9
    12 = 786
    13 = 27.7
    14 = 6.8
    16 = 26.9
    17 = 15.9
    18 = 21.7
    19 = 29.9
    20 = 30.2
    21 = 24
    22 = 31.6
    23 = 30.3
    24 = 20.7
    25 = 12
}

foreach ($row in $values.Keys) {
    $src = $ws.Range("N$row")
    $dst = $ws.Range("O$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $values[$row]
}

# Row 15 has no value in column N/O, only matching formatting.
$ws.Range("N15").Copy()
$ws.Range("O15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection, as recorded in the saved view state.
$ws.Range("Q20").Select()
